$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column A (rows 2-12) is treated as text so the date-like strings are not auto-converted to date serials
$ws.Range("A2:A12").NumberFormat = "@"

$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

$ws.Range("A2").Value = "2025-03-30"
$ws.Range("B2").Value = 0.3149557738767033
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("A3").Value = "2025-04-15"
$ws.Range("B3").Value = 0.2999517025687146
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.036136496927096855
$ws.Range("E3").Value = 0.0024186164830646982
$ws.Range("F3").Value = 0.009249891549247562
$ws.Range("G3").Value = 0.012021500634562681
$ws.Range("H3").Value = 0.002569579279478566
$ws.Range("I3").Value = -0.0035193048376580003
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -0.001607857489587361

$ws.Range("A4").Value = "2025-04-30"
$ws.Range("B4").Value = 0.3148344080686091
$ws.Range("C4").Value = 0.02628158220614346
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.015362852253286603
$ws.Range("F4").Value = -0.0014388995676532894
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.0009941533763434963
$ws.Range("I4").Value = -0.023996567393343906
$ws.Range("J4").Value = -0.00029823583450753826
$ws.Range("K4").Value = -0.002022179540374369

$ws.Range("A5").Value = "2025-05-15"
$ws.Range("B5").Value = 0.32957681152449203
$ws.Range("C5").Value = 0.09344954068926997
$ws.Range("D5").Value = -0.08582831948123892
$ws.Range("E5").Value = -0.002831419804499683
$ws.Range("F5").Value = 0.0494277101681891
$ws.Range("G5").Value = -0.04234004063183202
$ws.Range("H5").Value = 0.0015486764363057378
$ws.Range("I5").Value = 0.0015220017230588548
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.00020574564337005885

$ws.Range("A6").Value = "2025-05-30"
$ws.Range("B6").Value = 0.1940758272419768
$ws.Range("C6").Value = -0.0001487364395650905
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -0.1516242291503907
$ws.Range("F6").Value = 0.0009663720324204445
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.0068568891634963715
$ws.Range("I6").Value = 0.010365537395881377
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -0.0019168172843576081

$ws.Range("A7").Value = "2025-06-15"
$ws.Range("B7").Value = 0.15199332735154328
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.09030264296679949
$ws.Range("E7").Value = -0.021788883724172533
$ws.Range("F7").Value = -0.16301558838003802
$ws.Range("G7").Value = 0.009465659025884005
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.00641410782484339
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.03653956239625014

$ws.Range("A8").Value = "2025-06-30"
$ws.Range("B8").Value = -0.01273268742377659
$ws.Range("C8").Value = -0.15856457677863064
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.05881372144662604
$ws.Range("F8").Value = 0.005619014822953534
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = -0.00011303728995474196
$ws.Range("I8").Value = -0.06677614665505452
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = -0.003704990321259538

$ws.Range("A9").Value = "2025-07-15"
$ws.Range("B9").Value = -0.12383135141753154
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = -0.12920472074159292
$ws.Range("E9").Value = -0.015055349818621153
$ws.Range("F9").Value = 0.06488548510869441
$ws.Range("G9").Value = -0.01324981174296816
$ws.Range("H9").Value = 0.00034039143958367307
$ws.Range("I9").Value = -0.0016979732422327184
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = -0.017116684996618087

$ws.Range("A10").Value = "2025-07-30"
$ws.Range("B10").Value = 0.14340453251790483
$ws.Range("C10").Value = 0.35469750838914704
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -0.0019693849580461886
$ws.Range("F10").Value = -0.0034653549814477828
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = -0.0004411483192625705
$ws.Range("I10").Value = 0.012905574257154887
$ws.Range("J10").Value = -0.10390905555104726
$ws.Range("K10").Value = 0.009417745098938196

$ws.Range("A11").Value = "2025-08-15"
$ws.Range("B11").Value = 0.49299077886473064
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.2672178776126615
$ws.Range("E11").Value = 0.012315644453660223
$ws.Range("F11").Value = 0.022643520796621137
$ws.Range("G11").Value = -0.04121078472905382
$ws.Range("H11").Value = -0.003489887504376305
$ws.Range("I11").Value = 0.08205766798460243
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0.010052207732710622

$ws.Range("A12").Value = "2025-08-30"
$ws.Range("B12").Value = 0.11704734015663382
$ws.Range("C12").Value = -0.31449776328142265
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.000085196468178857175622409326
$ws.Range("F12").Value = -0.001276573722657786
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = -0.00015170472543905443
$ws.Range("I12").Value = -0.01454784024952116
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = -0.04555475319723501

# Reset column A style back to default (clears the custom text-format style) while keeping the text cell type
$ws.Range("A1:A12").Style = "Normal"

# Match the updated column widths (content-driven resize observed for columns C and F)
$ws.Columns.Item(3).ColumnWidth = 15.25
$ws.Columns.Item(6).ColumnWidth = 14.75